# Apply cryptos list update (Thu Sep  5 03:19:59 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.020.44"
$ws.Range("E2").Value = "  +0.71%  "

$ws.Range("D3").Value = "2.403.13"
$ws.Range("E3").Value = "  +1.84%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.06"
$ws.Range("E5").Value = "  -1.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.80"
$ws.Range("E6").Value = "  +3.12%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.553"
$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "2.417.02"
$ws.Range("E9").Value = "  +1.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0964"
$ws.Range("E10").Value = "  +0.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -0.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.321"
$ws.Range("E12").Value = "  +1.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.59"
$ws.Range("E13").Value = "  -4.22%  "

$ws.Range("D14").Value = "2.835.78"
$ws.Range("E14").Value = "  +2.04%  "

$ws.Range("D15").Value = "56.943.22"
$ws.Range("E15").Value = "  +0.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.74"
$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  +2.12%  "

$ws.Range("D18").Value = "2.438.99"
$ws.Range("E18").Value = "  +2.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.19"
$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "311.37"
$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.03"
$ws.Range("E21").Value = "  +0.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.34"
$ws.Range("E22").Value = "  +4.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.80"
$ws.Range("E23").Value = "  -0.73%  "

$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.33"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.154"
$ws.Range("E27").Value = "  +0.28%  "

$ws.Range("B28").Value = "Polygon"
$ws.Range("C28").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.379"
$ws.Range("E28").Value = "  -2.87%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.48"
$ws.Range("E29").Value = "  +4.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.90"
$ws.Range("E30").Value = "  -1.51%  "

$ws.Range("D31").Value = "0.0₃0722"
$ws.Range("E31").Value = "  +1.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.67"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.94"
$ws.Range("E33").Value = "  -2.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.11"
$ws.Range("E34").Value = "  -1.17%  "

$ws.Range("E35").Value = "  +0.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.91"
$ws.Range("E37").Value = "  +1.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.20"
$ws.Range("E38").Value = "  +1.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.83"
$ws.Range("E39").Value = "  +3.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.61"
$ws.Range("E40").Value = "  +3.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.800"
$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.43"
$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "130.82"
$ws.Range("E43").Value = "  +6.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.97"
$ws.Range("E44").Value = "  +1.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.36"
$ws.Range("E45").Value = "  +1.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "251.70"
$ws.Range("E46").Value = "  -0.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.564"
$ws.Range("E47").Value = "  -0.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0910"
$ws.Range("E48").Value = "  +0.41%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0487"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.00"
$ws.Range("E50").Value = "  +2.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0209"
$ws.Range("E51").Value = "  +1.17%  "
